$wb = $excel.ActiveWorkbook

# 1. Rename sheet "timepars_all" -> "timepars_old"
$wsOld = $wb.Worksheets.Item("timepars_all")
$wsOld.Name = "timepars_old"

# 2. On "timepars" sheet: add a comment at C1, change selection, and update I4/J4 values
$ws = $wb.Worksheets.Item("timepars")

$ws.Range("C1").AddComment("Effectiveness of PrEP")

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1

$ws.Range("F40").Select()
